$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: fill in the new task entry
# D18 - "2.5" needs to be stored as TEXT (like the existing "5.5" in D15), not as a number,
# while keeping the original cell style (s=2). Writing the string directly gets
# auto-coerced to a number, so instead we build it as a text formula and then
# paste back just the value, which keeps it as a shared string without changing style.
# (Done first so the new shared string for "2.5" is registered before the task text,
# matching the original authoring order.)
$ws.Range("D18").Formula = "=""2.5"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)

# B18 - task description (plain text, not numeric-looking so stays a normal string)
$ws.Range("B18").Value = "Ajout des dernières pages du site avec leurs laisons correspondantes, réorganisation des fichiers et renommage de certain. Amélioration du visuel du site"

# C18 - date (21/11/2015 = Excel serial 42329), reuse the existing date style (s=3)
$ws.Range("C18").Value = 42329
$ws.Range("C18").NumberFormat = "d-mmm"

# Update the active selection to B19, as recorded in the sheet view
$ws.Range("B19").Select()
